# Auto commit at 2025-12-27  8:07:35.66
# Append two new daily rows (2025-12-26 / serial 46017) for both stations
# to the bottom of the existing data table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Duplicate the formatting (date / decimal / integer cell styles) of the
# last existing data rows onto the two new rows, then overwrite the values.
$ws.Range("A50:F51").Copy()
$ws.Range("A52").PasteSpecial()

# Row 52: 四方坪站 (station 1)
$ws.Range("A52").Value = 46017
$ws.Range("B52").Value = "四方坪站"
$ws.Range("C52").Value = 10476.65
$ws.Range("D52").Value = 9057.5499999999993
$ws.Range("E52").Value = 3453.99
$ws.Range("F52").Value = 423

# Row 53: 高岭站 (station 2)
$ws.Range("A53").Value = 46017
$ws.Range("B53").Value = "高岭站"
$ws.Range("C53").Value = 6068.09
$ws.Range("D53").Value = 5177.62
$ws.Range("E53").Value = 1619.76
$ws.Range("F53").Value = 203

# Update selection to match the target state (cursor moved to K50 after edits)
$ws.Range("K50").Select()
